# Source data format changed from Excel to CSV, and the "chiclet" output
# column is split into separate "chiclet-version" / "chiclet-PI" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D ("fireTV" and everything after it
# shifts one column to the right: D->E, E->F, F->G, G->H, H->I, I->J).
[void]$ws.Columns.Item(4).Insert()

# New column D becomes "chiclet-PI" with a numeric PI/date value.
$ws.Range("D1").Value2 = "chiclet-PI"
$ws.Range("D2").Value2 = 20210224
$ws.Columns.Item(4).ColumnWidth = 11.83

# Existing column C ("chiclet") becomes "chiclet-version" with the new
# version-string value (instead of the old numeric year value).
$ws.Range("C1").Value2 = "chiclet-version"
$ws.Range("C2").Value2 = "8.9.0.36063"

# The old "Playback Hours" column (shifted from G to H by the insert above)
# is removed entirely - chromecast/tvOS slide back into H/I.
[void]$ws.Columns.Item(8).Delete()

# Update the active selection to match the new layout.
[void]$ws.Range("H2").Select()
